# Edit script: add "Geocomputation with R" resource (two rows), add three new
# authors (Robin Lovelace, Jakub Nowosad, Jannes Muenchow) to the Authors sheet.

$wb = $excel.ActiveWorkbook
$wsResources = $wb.Worksheets.Item("Resources")
$wsAuthors = $wb.Worksheets.Item("Authors")

# ---------------------------------------------------------------------------
# 1. Resources sheet: insert two new rows right after the header row and fill
#    them in with the new "Geocomputation with R" entries (one tagged
#    Programming, one tagged Visualisation).
# ---------------------------------------------------------------------------
$wsResources.Rows("2:3").Insert()

# The insert operation can drag formatting/hyperlinks down from row 1 (or
# duplicate row 2's old hyperlink) onto the freshly inserted rows - wipe them
# back to a clean slate before writing the new values.
$wsResources.Rows("2:3").Hyperlinks.Delete()
$wsResources.Range("A2:I3").ClearContents()
$wsResources.Range("A2:I3").ClearFormats()

$wsResources.Range("A2").Value = "Web"
$wsResources.Range("B2").Value = "Programming"
$wsResources.Range("C2").Value = "Reference"
$wsResources.Range("D2").Value = "R"
$wsResources.Range("F2").Value = "Geocomputation with R"
$wsResources.Range("G2").Value = "Robin Lovelace, Jakub Nowosad, Jannes Muenchow"
$wsResources.Range("H2").Value = "https://geocompr.robinlovelace.net"
$wsResources.Hyperlinks.Add($wsResources.Range("H2"), "https://geocompr.robinlovelace.net") | Out-Null

$wsResources.Range("A3").Value = "Web"
$wsResources.Range("B3").Value = "Visualisation"
$wsResources.Range("C3").Value = "Reference"
$wsResources.Range("D3").Value = "R"
$wsResources.Range("F3").Value = "Geocomputation with R"
$wsResources.Range("G3").Value = "Robin Lovelace, Jakub Nowosad, Jannes Muenchow"
$wsResources.Range("H3").Value = "https://geocompr.robinlovelace.net"
$wsResources.Hyperlinks.Add($wsResources.Range("H3"), "https://geocompr.robinlovelace.net") | Out-Null

# ---------------------------------------------------------------------------
# 2. Authors sheet: append the three new authors related to the book above.
# ---------------------------------------------------------------------------
$wsAuthors.Range("A35").Value = "Robin Lovelace"
$wsAuthors.Range("B35").Value = "https://www.robinlovelace.net"

$wsAuthors.Range("A36").Value = "Jakub Nowosad"
$wsAuthors.Range("B36").Value = "https://nowosad.github.io"

$wsAuthors.Range("A37").Value = "Jannes Muenchow"
$wsAuthors.Range("B37").Value = "https://www.geographie.uni-jena.de/en/Muenchow.html"
$wsAuthors.Hyperlinks.Add($wsAuthors.Range("B37"), "https://www.geographie.uni-jena.de/en/Muenchow.html") | Out-Null

# ---------------------------------------------------------------------------
# 3. Restore sensible view state: Authors sheet scrolled/selected at B38 (the
#    first empty row below the new authors), then flip back to Resources with
#    C3:H3 selected so it remains the active sheet/tab on save.
# ---------------------------------------------------------------------------
$wsAuthors.Select()
$wsAuthors.Range("B38").Select()

$wsResources.Select()
$wsResources.Range("C3:H3").Select()

Write-Output "done"
